$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "list": insert 4 new columns (TA, UM, ND, AZ) right before
# the existing "ALL" column (column I), pushing ALL to column M.
# ---------------------------------------------------------------
$wsList = $wb.Worksheets.Item("list")

$wsList.Columns.Item(9).Insert()
$wsList.Columns.Item(9).Insert()
$wsList.Columns.Item(9).Insert()
$wsList.Columns.Item(9).Insert()

$wsList.Range("I1").Value = "TA"
$wsList.Range("J1").Value = "UM"
$wsList.Range("K1").Value = "ND"
$wsList.Range("L1").Value = "AZ"

for ($r = 2; $r -le 16; $r++) {
    $wsList.Cells.Item($r, 9).Value = 0
    $wsList.Cells.Item($r, 10).Value = 0
    $wsList.Cells.Item($r, 11).Value = 0
    $wsList.Cells.Item($r, 12).Value = 0
}

# ---------------------------------------------------------------
# Sheet "summary": insert 4 new rows (TA, UM, ND, AZ) right before
# the existing "ALL" row (row 7), pushing ALL to row 11.
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("summary")

$wsSummary.Rows.Item(7).Insert()
$wsSummary.Rows.Item(7).Insert()
$wsSummary.Rows.Item(7).Insert()
$wsSummary.Rows.Item(7).Insert()

$labels = @("TA", "UM", "ND", "AZ")
for ($i = 0; $i -lt 4; $i++) {
    $r = 7 + $i
    $wsSummary.Cells.Item($r, 1).Value = $labels[$i]
    $wsSummary.Cells.Item($r, 2).Value = 0
    $wsSummary.Cells.Item($r, 3).Value = 0
    $wsSummary.Cells.Item($r, 4).Value = 0
    $wsSummary.Cells.Item($r, 5).Value = 0
    $wsSummary.Cells.Item($r, 6).Value = 0
}
